$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticker data set: replaces the old A2:C4 rows with an expanded A2:C8 table
$data = @(
    @(1539, "KB",   0),
    @(1744, "CIB",  0),
    @(2211, "AVAL", 0),
    @(3608, "JHB",  0),
    @(4095, "DMTK", 1),
    @(4848, "GLBD", 0),
    @(4995, "MMNT", 0)
)

# Column A uses the same bordered / bold / centered formatting throughout;
# grab it once from the still-formatted A2 before any values change.
$formatSource = $ws.Range("A2")
$formatSource.Copy() | Out-Null

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Re-apply / extend column A's style onto the newly added rows (5-8)
$ws.Range("A5:A8").PasteSpecial(-4122) | Out-Null
